$wb = $excel.ActiveWorkbook

# ===== Sheet1: Chateau =====
$ws1 = $wb.Worksheets.Item("Chateau")
$ws1.Range("C77").Value = 1735
$ws1.Range("D77").Value = 26
$ws1.Range("E77").Value = 14
$ws1.Range("C78").Value = 795
$ws1.Range("D78").Value = 14
$ws1.Range("E78").Value = 11
$ws1.Range("C79").Value = 540
$ws1.Range("D79").Value = 10
$ws1.Range("E79").Value = 8
$ws1.Range("C80").Value = 445
$ws1.Range("D80").Value = 13
$ws1.Range("E80").Value = 17
$ws1.Range("C81").Value = 465
$ws1.Range("D81").Value = 10
$ws1.Range("E81").Value = 13
$ws1.Range("A81:G81").Copy($ws1.Range("A82:G82"))
$ws1.Range("A82").Value = 10
$ws1.Range("B82").Value = "Juanito"
$ws1.Range("C82").Value = 0
$ws1.Range("D82").Value = 10
$ws1.Range("E82").Value = 21
$ws1.Range("F82").Formula = "=IF(E82=0,0,(D82/E82))"
$ws1.Range("A82:G82").Copy($ws1.Range("A83:G83"))
$ws1.Range("A83").Value = 11
$ws1.Range("B83").Value = "Fer"
$ws1.Range("C83").Value = 1185
$ws1.Range("D83").Value = 22
$ws1.Range("E83").Value = 11
$ws1.Range("F83").Formula = "=IF(E83=0,0,(D83/E83))"
$ws1.Range("A83:G83").Copy($ws1.Range("A84:G84"))
$ws1.Range("A84").Value = 11
$ws1.Range("B84").Value = "Joako"
$ws1.Range("C84").Value = 755
$ws1.Range("D84").Value = 13
$ws1.Range("E84").Value = 10
$ws1.Range("F84").Formula = "=IF(E84=0,0,(D84/E84))"
$ws1.Range("A84:G84").Copy($ws1.Range("A85:G85"))
$ws1.Range("A85").Value = 11
$ws1.Range("B85").Value = "Arturo"
$ws1.Range("C85").Value = 645
$ws1.Range("D85").Value = 13
$ws1.Range("E85").Value = 14
$ws1.Range("F85").Formula = "=IF(E85=0,0,(D85/E85))"
$ws1.Range("A85:G85").Copy($ws1.Range("A86:G86"))
$ws1.Range("A86").Value = 11
$ws1.Range("B86").Value = "Oscar"
$ws1.Range("C86").Value = 335
$ws1.Range("D86").Value = 13
$ws1.Range("E86").Value = 15
$ws1.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws1.Range("A86:G86").Copy($ws1.Range("A87:G87"))
$ws1.Range("A87").Value = 11
$ws1.Range("B87").Value = "Brandon"
$ws1.Range("C87").Value = 585
$ws1.Range("D87").Value = 16
$ws1.Range("E87").Value = 11
$ws1.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws1.Range("A87:G87").Copy($ws1.Range("A88:G88"))
$ws1.Range("A88").Value = 11
$ws1.Range("B88").Value = "Edson"
$ws1.Range("C88").Value = 275
$ws1.Range("D88").Value = 3
$ws1.Range("E88").Value = 9
$ws1.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws1.Range("A88:G88").Copy($ws1.Range("A89:G89"))
$ws1.Range("A89").Value = 11
$ws1.Range("B89").Value = "Juanito"
$ws1.Range("C89").Value = 10
$ws1.Range("D89").Value = 7
$ws1.Range("E89").Value = 17
$ws1.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"

# ===== Sheet2: Centro Alpino =====
$ws2 = $wb.Worksheets.Item("Centro Alpino")
$ws2.Range("A85:G85").Copy($ws2.Range("A86:G86"))
$ws2.Range("A86").Value = 11
$ws2.Range("B86").Value = "Fer"
$ws2.Range("C86").Value = 490
$ws2.Range("D86").Value = 8
$ws2.Range("E86").Value = 6
$ws2.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws2.Range("A86:G86").Copy($ws2.Range("A87:G87"))
$ws2.Range("A87").Value = 11
$ws2.Range("B87").Value = "Joako"
$ws2.Range("C87").Value = 140
$ws2.Range("D87").Value = 4
$ws2.Range("E87").Value = 5
$ws2.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws2.Range("A87:G87").Copy($ws2.Range("A88:G88"))
$ws2.Range("A88").Value = 11
$ws2.Range("B88").Value = "Arturo"
$ws2.Range("C88").Value = 585
$ws2.Range("D88").Value = 6
$ws2.Range("E88").Value = 6
$ws2.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws2.Range("A88:G88").Copy($ws2.Range("A89:G89"))
$ws2.Range("A89").Value = 11
$ws2.Range("B89").Value = "Oscar"
$ws2.Range("C89").Value = 1740
$ws2.Range("D89").Value = 8
$ws2.Range("E89").Value = 0
$ws2.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"
$ws2.Range("A89:G89").Copy($ws2.Range("A90:G90"))
$ws2.Range("A90").Value = 11
$ws2.Range("B90").Value = "Brandon"
$ws2.Range("C90").Value = 0
$ws2.Range("D90").Value = 3
$ws2.Range("E90").Value = 8
$ws2.Range("F90").Formula = "=IF(E90=0,0,(D90/E90))"
$ws2.Range("A90:G90").Copy($ws2.Range("A91:G91"))
$ws2.Range("A91").Value = 11
$ws2.Range("B91").Value = "Edson"
$ws2.Range("C91").Value = 45
$ws2.Range("D91").Value = 4
$ws2.Range("E91").Value = 5
$ws2.Range("F91").Formula = "=IF(E91=0,0,(D91/E91))"
$ws2.Range("A91:G91").Copy($ws2.Range("A92:G92"))
$ws2.Range("A92").Value = 11
$ws2.Range("B92").Value = "Juanito"
$ws2.Range("C92").Value = 55
$ws2.Range("D92").Value = 4
$ws2.Range("E92").Value = 5
$ws2.Range("F92").Formula = "=IF(E92=0,0,(D92/E92))"

# ===== Sheet3: Ciudad Destruida =====
$ws3 = $wb.Worksheets.Item("Ciudad Destruida")
$ws3.Range("A84:G84").Copy($ws3.Range("A85:G85"))
$ws3.Range("A85").Value = 11
$ws3.Range("B85").Value = "Fer"
$ws3.Range("C85").Value = 1375
$ws3.Range("D85").Value = 16
$ws3.Range("E85").Value = 8
$ws3.Range("F85").Formula = "=IF(E85=0,0,(D85/E85))"
$ws3.Range("A85:G85").Copy($ws3.Range("A86:G86"))
$ws3.Range("A86").Value = 11
$ws3.Range("B86").Value = "Joako"
$ws3.Range("C86").Value = 100
$ws3.Range("D86").Value = 5
$ws3.Range("E86").Value = 12
$ws3.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws3.Range("A86:G86").Copy($ws3.Range("A87:G87"))
$ws3.Range("A87").Value = 11
$ws3.Range("B87").Value = "Arturo"
$ws3.Range("C87").Value = 10
$ws3.Range("D87").Value = 9
$ws3.Range("E87").Value = 10
$ws3.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws3.Range("A87:G87").Copy($ws3.Range("A88:G88"))
$ws3.Range("A88").Value = 11
$ws3.Range("B88").Value = "Oscar"
$ws3.Range("C88").Value = 635
$ws3.Range("D88").Value = 9
$ws3.Range("E88").Value = 8
$ws3.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws3.Range("A88:G88").Copy($ws3.Range("A89:G89"))
$ws3.Range("A89").Value = 11
$ws3.Range("B89").Value = "Brandon"
$ws3.Range("C89").Value = 525
$ws3.Range("D89").Value = 10
$ws3.Range("E89").Value = 7
$ws3.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"
$ws3.Range("A89:G89").Copy($ws3.Range("A90:G90"))
$ws3.Range("A90").Value = 11
$ws3.Range("B90").Value = "Edson"
$ws3.Range("C90").Value = 645
$ws3.Range("D90").Value = 4
$ws3.Range("E90").Value = 6
$ws3.Range("F90").Formula = "=IF(E90=0,0,(D90/E90))"
$ws3.Range("A90:G90").Copy($ws3.Range("A91:G91"))
$ws3.Range("A91").Value = 11
$ws3.Range("B91").Value = "Juanito"
$ws3.Range("C91").Value = 260
$ws3.Range("D91").Value = 6
$ws3.Range("E91").Value = 9
$ws3.Range("F91").Formula = "=IF(E91=0,0,(D91/E91))"

# ===== Sheet4: Patio Maniobras =====
$ws4 = $wb.Worksheets.Item("Patio Maniobras")
$ws4.Range("A85:G85").Copy($ws4.Range("A86:G86"))
$ws4.Range("A86").Value = 11
$ws4.Range("B86").Value = "Fer"
$ws4.Range("C86").Value = 1335
$ws4.Range("D86").Value = 17
$ws4.Range("E86").Value = 10
$ws4.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws4.Range("A86:G86").Copy($ws4.Range("A87:G87"))
$ws4.Range("A87").Value = 11
$ws4.Range("B87").Value = "Joako"
$ws4.Range("C87").Value = 1285
$ws4.Range("D87").Value = 13
$ws4.Range("E87").Value = 13
$ws4.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws4.Range("A87:G87").Copy($ws4.Range("A88:G88"))
$ws4.Range("A88").Value = 11
$ws4.Range("B88").Value = "Arturo"
$ws4.Range("C88").Value = 715
$ws4.Range("D88").Value = 15
$ws4.Range("E88").Value = 17
$ws4.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws4.Range("A88:G88").Copy($ws4.Range("A89:G89"))
$ws4.Range("A89").Value = 11
$ws4.Range("B89").Value = "Oscar"
$ws4.Range("C89").Value = 425
$ws4.Range("D89").Value = 12
$ws4.Range("E89").Value = 10
$ws4.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"
$ws4.Range("A89:G89").Copy($ws4.Range("A90:G90"))
$ws4.Range("A90").Value = 11
$ws4.Range("B90").Value = "Brandon"
$ws4.Range("C90").Value = 885
$ws4.Range("D90").Value = 14
$ws4.Range("E90").Value = 8
$ws4.Range("F90").Formula = "=IF(E90=0,0,(D90/E90))"
$ws4.Range("A90:G90").Copy($ws4.Range("A91:G91"))
$ws4.Range("A91").Value = 11
$ws4.Range("B91").Value = "Edson"
$ws4.Range("C91").Value = 0
$ws4.Range("D91").Value = 4
$ws4.Range("E91").Value = 12
$ws4.Range("F91").Formula = "=IF(E91=0,0,(D91/E91))"
$ws4.Range("A91:G91").Copy($ws4.Range("A92:G92"))
$ws4.Range("A92").Value = 11
$ws4.Range("B92").Value = "Juanito"
$ws4.Range("C92").Value = 375
$ws4.Range("D92").Value = 6
$ws4.Range("E92").Value = 11
$ws4.Range("F92").Formula = "=IF(E92=0,0,(D92/E92))"

# ===== Sheet5: Saint Raymonds =====
$ws5 = $wb.Worksheets.Item("Saint Raymonds")
$ws5.Range("A84:G84").Copy($ws5.Range("A85:G85"))
$ws5.Range("A85").Value = 11
$ws5.Range("B85").Value = "Fer"
$ws5.Range("C85").Value = 615
$ws5.Range("D85").Value = 23
$ws5.Range("E85").Value = 16
$ws5.Range("F85").Formula = "=IF(E85=0,0,(D85/E85))"
$ws5.Range("A85:G85").Copy($ws5.Range("A86:G86"))
$ws5.Range("A86").Value = 11
$ws5.Range("B86").Value = "Joako"
$ws5.Range("C86").Value = 1740
$ws5.Range("D86").Value = 17
$ws5.Range("E86").Value = 10
$ws5.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws5.Range("A86:G86").Copy($ws5.Range("A87:G87"))
$ws5.Range("A87").Value = 11
$ws5.Range("B87").Value = "Arturo"
$ws5.Range("C87").Value = 0
$ws5.Range("D87").Value = 9
$ws5.Range("E87").Value = 22
$ws5.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws5.Range("A87:G87").Copy($ws5.Range("A88:G88"))
$ws5.Range("A88").Value = 11
$ws5.Range("B88").Value = "Oscar"
$ws5.Range("C88").Value = 955
$ws5.Range("D88").Value = 14
$ws5.Range("E88").Value = 7
$ws5.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws5.Range("A88:G88").Copy($ws5.Range("A89:G89"))
$ws5.Range("A89").Value = 11
$ws5.Range("B89").Value = "Brandon"
$ws5.Range("C89").Value = 0
$ws5.Range("D89").Value = 6
$ws5.Range("E89").Value = 12
$ws5.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"
$ws5.Range("A89:G89").Copy($ws5.Range("A90:G90"))
$ws5.Range("A90").Value = 11
$ws5.Range("B90").Value = "Edson"
$ws5.Range("C90").Value = 375
$ws5.Range("D90").Value = 4
$ws5.Range("E90").Value = 6
$ws5.Range("F90").Formula = "=IF(E90=0,0,(D90/E90))"
$ws5.Range("A90:G90").Copy($ws5.Range("A91:G91"))
$ws5.Range("A91").Value = 11
$ws5.Range("B91").Value = "Juanito"
$ws5.Range("C91").Value = 215
$ws5.Range("D91").Value = 15
$ws5.Range("E91").Value = 16
$ws5.Range("F91").Formula = "=IF(E91=0,0,(D91/E91))"

# ===== Sheet6: Zona Aterrizaje =====
$ws6 = $wb.Worksheets.Item("Zona Aterrizaje")
$ws6.Range("A83:G83").Copy($ws6.Range("A84:G84"))
$ws6.Range("A84").Value = 10
$ws6.Range("B84").Value = "Juanito"
$ws6.Range("C84").Value = 105
$ws6.Range("D84").Value = 2
$ws6.Range("E84").Value = 12
$ws6.Range("F84").Formula = "=IF(E84=0,0,(D84/E84))"
$ws6.Range("A84:G84").Copy($ws6.Range("A85:G85"))
$ws6.Range("A85").Value = 11
$ws6.Range("B85").Value = "Fer"
$ws6.Range("C85").Value = 1130
$ws6.Range("D85").Value = 18
$ws6.Range("E85").Value = 15
$ws6.Range("F85").Formula = "=IF(E85=0,0,(D85/E85))"
$ws6.Range("A85:G85").Copy($ws6.Range("A86:G86"))
$ws6.Range("A86").Value = 11
$ws6.Range("B86").Value = "Joako"
$ws6.Range("C86").Value = 1145
$ws6.Range("D86").Value = 13
$ws6.Range("E86").Value = 8
$ws6.Range("F86").Formula = "=IF(E86=0,0,(D86/E86))"
$ws6.Range("A86:G86").Copy($ws6.Range("A87:G87"))
$ws6.Range("A87").Value = 11
$ws6.Range("B87").Value = "Arturo"
$ws6.Range("C87").Value = 240
$ws6.Range("D87").Value = 14
$ws6.Range("E87").Value = 16
$ws6.Range("F87").Formula = "=IF(E87=0,0,(D87/E87))"
$ws6.Range("A87:G87").Copy($ws6.Range("A88:G88"))
$ws6.Range("A88").Value = 11
$ws6.Range("B88").Value = "Oscar"
$ws6.Range("C88").Value = 1040
$ws6.Range("D88").Value = 20
$ws6.Range("E88").Value = 12
$ws6.Range("F88").Formula = "=IF(E88=0,0,(D88/E88))"
$ws6.Range("A88:G88").Copy($ws6.Range("A89:G89"))
$ws6.Range("A89").Value = 11
$ws6.Range("B89").Value = "Brandon"
$ws6.Range("C89").Value = 150
$ws6.Range("D89").Value = 7
$ws6.Range("E89").Value = 16
$ws6.Range("F89").Formula = "=IF(E89=0,0,(D89/E89))"
$ws6.Range("A89:G89").Copy($ws6.Range("A90:G90"))
$ws6.Range("A90").Value = 11
$ws6.Range("B90").Value = "Edson"
$ws6.Range("C90").Value = 795
$ws6.Range("D90").Value = 7
$ws6.Range("E90").Value = 9
$ws6.Range("F90").Formula = "=IF(E90=0,0,(D90/E90))"
$ws6.Range("A90:G90").Copy($ws6.Range("A91:G91"))
$ws6.Range("A91").Value = 11
$ws6.Range("B91").Value = "Juanito"
$ws6.Range("C91").Value = 545
$ws6.Range("D91").Value = 5
$ws6.Range("E91").Value = 10
$ws6.Range("F91").Formula = "=IF(E91=0,0,(D91/E91))"

# ===== Restore on-screen selections to match the saved view state =====
$ws1.Range("B83:B89").Select()
$ws2.Range("C93").Select()
$ws3.Range("F91").Select()
$ws5.Range("C92").Select()
$ws6.Range("F94").Select()
$ws4.Activate()
$ws4.Range("F94").Select()

Write-Host "Edit complete"
